$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7176
$ws.Range("C3").Value = 157178
$ws.Range("C4").Value = 148244
$ws.Range("C5").Value = 8934
$ws.Range("C7").Value = 5.68
$ws.Range("C8").Value = 63.84
